# Commit "#5: pham vi du an"
#
# The ToC bookmark "_Toc57272150" (which marks the "Giao tiep/Trao doi
# thong tin" Heading 1) is relocated so its *start* sits at the top of
# the preceding paragraph instead -- the short italic blurb right
# before that heading. That paragraph's wording is also rewritten
# ("Chay duoc tren nen tang nao, OS nao?..." -> "Chay tren he dieu
# hanh Window") and its italic formatting is removed; the heading
# paragraph itself keeps the bookmark's end mark but loses the
# bookmark's start mark and is otherwise unchanged.

$d = $word.ActiveDocument

# Anchor on the bookmark itself (sturdier than matching the old
# Vietnamese text, and it's exactly the element the edit relocates).
$bookmark = $d.Bookmarks("_Toc57272150")
$headingPara = $bookmark.Range.Paragraphs(1)
$blurbPara = $headingPara.Previous()

$editRange = $d.Range($blurbPara.Range.Start, $headingPara.Range.End)
$editRange.InsertXML('<w:p w14:paraId="303DB423" w14:textId="7BAE1263" w:rsidR="00A01B4E" w:rsidRPr="00A01B4E" w:rsidRDefault="00A01B4E" w:rsidP="00A01B4E"><w:bookmarkStart w:id="12" w:name="_Toc57272150"/><w:proofErr w:type="spellStart"/><w:r><w:t>Chạy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trên</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hệ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>điều</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hành</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Window</w:t></w:r></w:p><w:p w14:paraId="51906836" w14:textId="543D8E3B" w:rsidR="00802E21" w:rsidRDefault="00802E21" w:rsidP="00F16A81"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">Giao </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tiếp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Trao</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>đổi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thông</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tin</w:t></w:r><w:bookmarkEnd w:id="12"/></w:p>')
